$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "89÷4=22, 1"
$tbl.Cell(1,2).Range.Text = "76÷5=15, 1"
$tbl.Cell(1,3).Range.Text = "97÷7=13, 6"
$tbl.Cell(1,4).Range.Text = "29÷5=5, 4"
$tbl.Cell(1,5).Range.Text = "65÷4=16, 1"

$tbl.Cell(5,1).Range.Text = "67÷9=7, 4"
$tbl.Cell(5,2).Range.Text = "18÷3=6, 0"
$tbl.Cell(5,3).Range.Text = "17÷7=2, 3"
$tbl.Cell(5,4).Range.Text = "41÷4=10, 1"
$tbl.Cell(5,5).Range.Text = "41÷4=10, 1"

$tbl.Cell(9,1).Range.Text = "96÷4=24, 0"
$tbl.Cell(9,2).Range.Text = "14÷6=2, 2"
$tbl.Cell(9,3).Range.Text = "17÷9=1, 8"
$tbl.Cell(9,4).Range.Text = "87÷2=43, 1"
$tbl.Cell(9,5).Range.Text = "32÷4=8, 0"

$tbl.Cell(13,1).Range.Text = "65÷5=13, 0"
$tbl.Cell(13,2).Range.Text = "86÷6=14, 2"
$tbl.Cell(13,3).Range.Text = "71÷4=17, 3"
$tbl.Cell(13,4).Range.Text = "41÷7=5, 6"
$tbl.Cell(13,5).Range.Text = "48÷4=12, 0"

$tbl.Cell(17,1).Range.Text = "91÷8=11, 3"
$tbl.Cell(17,2).Range.Text = "62÷5=12, 2"
$tbl.Cell(17,3).Range.Text = "63÷5=12, 3"
$tbl.Cell(17,4).Range.Text = "46÷5=9, 1"
$tbl.Cell(17,5).Range.Text = "95÷4=23, 3"
